# Latest update of prototyping
#
# 1. "inv_cost" gets new (recomputed) cost figures.
# 2. A brand-new "fix_cost" sheet is inserted (node_loc/technology/year_vtg/
#    year_act/value/unit) between "initial_new_capacity_up" and
#    "technical_lifetime".
# 3. Tabs are reordered to:
#    inv_cost, initial_new_capacity_up, fix_cost, technical_lifetime,
#    growth_new_capacity_up
# ("initial_new_capacity_up", "technical_lifetime" and
#  "growth_new_capacity_up" keep their existing data - only their tab
#  position changes.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update "inv_cost" values (column D) with the recomputed figures.
# ---------------------------------------------------------------------------
$invCost = $wb.Worksheets.Item("inv_cost")

$invCostValues = @(
    2500,
    1496.842348095947,
    896.2148060213547,
    2000,
    1197.473878476757,
    716.9718448170838
)

for ($i = 0; $i -lt $invCostValues.Length; $i++) {
    $row = $i + 2
    $invCost.Cells.Item($row, 4).Value2 = $invCostValues[$i]
}

# ---------------------------------------------------------------------------
# 2. Reorder existing tabs so only "fix_cost" needs to be inserted to reach
#    the final order:
#       inv_cost, initial_new_capacity_up, fix_cost, technical_lifetime,
#       growth_new_capacity_up
#
# NOTE: worksheet handles returned by Worksheets.Item(...) are positional
# snapshots that get re-resolved by their original index after the tab order
# changes - they are NOT stable identities. So every handle is re-fetched by
# name immediately before it is used, right after any Move()/Add() call.
# ---------------------------------------------------------------------------

# Move "technical_lifetime" after "growth_new_capacity_up" (to the very end).
$src1 = $wb.Worksheets.Item("technical_lifetime")
$dst1 = $wb.Worksheets.Item("growth_new_capacity_up")
$src1.Move($null, $dst1)

# Move "growth_new_capacity_up" after "technical_lifetime" (back to the end,
# after the sheet that just passed it). Re-fetch both handles: the previous
# Move() invalidated every previously-held worksheet reference.
$src2 = $wb.Worksheets.Item("growth_new_capacity_up")
$dst2 = $wb.Worksheets.Item("technical_lifetime")
$src2.Move($null, $dst2)

# Current order now: inv_cost, initial_new_capacity_up, technical_lifetime,
# growth_new_capacity_up

# ---------------------------------------------------------------------------
# 3. Insert the new "fix_cost" sheet right after "initial_new_capacity_up"
#    (i.e. right before "technical_lifetime"). Re-fetch again: the Move()
#    above invalidated old handles.
# ---------------------------------------------------------------------------
$initialNewCapacityUp = $wb.Worksheets.Item("initial_new_capacity_up")
$fixCost = $wb.Worksheets.Add($null, $initialNewCapacityUp)
$fixCost.Name = "fix_cost"

# Header row. Cells are written one at a time - assigning a (multi-cell)
# array straight to a Range's .Value2 silently produces no cell data in this
# host, unlike per-cell assignment.
$fixCost.Cells.Item(1, 1).Value2 = "node_loc"
$fixCost.Cells.Item(1, 2).Value2 = "technology"
$fixCost.Cells.Item(1, 3).Value2 = "year_vtg"
$fixCost.Cells.Item(1, 4).Value2 = "year_act"
$fixCost.Cells.Item(1, 5).Value2 = "value"
$fixCost.Cells.Item(1, 6).Value2 = "unit"

# Style the header row like the other sheets' header row (bold, thin border,
# centered horizontally, top-aligned vertically).
$header = $fixCost.Range("A1:F1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows: every node_loc x technology x year_vtg x year_act combination
# where year_act >= year_vtg.
$nodes = @("Westeros", "Westerlands")
$years = @(700, 710, 720)

$r = 1
foreach ($node in $nodes) {
    foreach ($vtg in $years) {
        foreach ($act in $years) {
            if ($act -ge $vtg) {
                $r = $r + 1
                $fixCost.Cells.Item($r, 1).Value2 = $node
                $fixCost.Cells.Item($r, 2).Value2 = "DACCS"
                $fixCost.Cells.Item($r, 3).Value2 = $vtg
                $fixCost.Cells.Item($r, 4).Value2 = $act
                $fixCost.Cells.Item($r, 5).Value2 = 10
                $fixCost.Cells.Item($r, 6).Value2 = "$/kWa"
            }
        }
    }
}

Write-Output "done"
